# This script updates the per-base-position frequency table (rows 2-5 = A/C/G/T,
# columns B:X = positions 1-23) with the re-run values reported for publication.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("A" base) - columns B:X
$ws.Cells.Item(2, 2).Value = 0.0995906397018391  # B2
$ws.Cells.Item(2, 3).Value = 0.968106555874626  # C2
$ws.Cells.Item(2, 4).Value = 0.000488788415714548  # D2
$ws.Cells.Item(2, 5).Value = 0.993890144803568  # E2
$ws.Cells.Item(2, 6).Value = 0.00837050161911163  # F2
$ws.Cells.Item(2, 7).Value = 0.955398057066048  # G2
$ws.Cells.Item(2, 8).Value = 0.0231563511944767  # H2
$ws.Cells.Item(2, 9).Value = 0.276348750534612  # I2
$ws.Cells.Item(2, 10).Value = 0.000244394207857274  # J2
$ws.Cells.Item(2, 11).Value = 0.95246532657176  # K2
$ws.Cells.Item(2, 12).Value = 0.991996089692674  # L2
$ws.Cells.Item(2, 13).Value = 0.874198081505468  # M2
$ws.Cells.Item(2, 14).Value = 0.943606036536934  # N2
$ws.Cells.Item(2, 15).Value = 0.000855379727500458  # O2
$ws.Cells.Item(2, 16).Value = 0.999205718824464  # P2
$ws.Cells.Item(2, 17).Value = 0.0000610985519643184  # Q2
$ws.Cells.Item(2, 18).Value = 0.999022423168571  # R2
$ws.Cells.Item(2, 19).Value = 0.0000610985519643184  # S2
$ws.Cells.Item(2, 20).Value = 0.00403250442964502  # T2
$ws.Cells.Item(2, 21).Value = 0.995295411498747  # U2
$ws.Cells.Item(2, 22).Value = 0.988024683814994  # V2
$ws.Cells.Item(2, 23).Value = 0.0887761960041547  # W2
$ws.Cells.Item(2, 24).Value = 0.000733182623571821  # X2

# Row 3 ("C" base) - columns B:X
$ws.Cells.Item(3, 2).Value = 0.0530946416569927  # B3
$ws.Cells.Item(3, 3).Value = 0.000305492759821592  # C3
$ws.Cells.Item(3, 4).Value = 0.00122197103928637  # D3
$ws.Cells.Item(3, 5).Value = 0.00128306959125069  # E3
$ws.Cells.Item(3, 6).Value = 0.00775951609946844  # F3
$ws.Cells.Item(3, 7).Value = 0.0152135394391153  # G3
$ws.Cells.Item(3, 8).Value = 0.0389197776012708  # H3
$ws.Cells.Item(3, 9).Value = 0.0000610985519643184  # I3
$ws.Cells.Item(3, 10).Value = 0.000549886967678866  # J3
$ws.Cells.Item(3, 11).Value = 0.000916478279464777  # K3
$ws.Cells.Item(3, 12).Value = 0.00378811022178774  # L3
$ws.Cells.Item(3, 13).Value = 0.00409360298160934  # M3
$ws.Cells.Item(3, 14).Value = 0.000610985519643184  # N3
$ws.Cells.Item(3, 15).Value = 0  # O3
$ws.Cells.Item(3, 16).Value = 0.0000610985519643184  # P3
$ws.Cells.Item(3, 17).Value = 0.917883546159956  # Q3
$ws.Cells.Item(3, 18).Value = 0.0000610985519643184  # R3
$ws.Cells.Item(3, 19).Value = 0.999450113032321  # S3
$ws.Cells.Item(3, 20).Value = 0.985702938840349  # T3
$ws.Cells.Item(3, 21).Value = 0.00342151891000183  # U3
$ws.Cells.Item(3, 22).Value = 0.000366591311785911  # V3
$ws.Cells.Item(3, 23).Value = 0.000183295655892955  # W3
$ws.Cells.Item(3, 24).Value = 0.000122197103928637  # X3

# Row 4 ("G" base) - columns B:X
$ws.Cells.Item(4, 2).Value = 0.774668540355594  # B4
$ws.Cells.Item(4, 3).Value = 0.027433249831979  # C4
$ws.Cells.Item(4, 4).Value = 0.000549886967678866  # D4
$ws.Cells.Item(4, 5).Value = 0.00354371601393047  # E4
$ws.Cells.Item(4, 6).Value = 0.982709109794098  # F4
$ws.Cells.Item(4, 7).Value = 0.00879819148286186  # G4
$ws.Cells.Item(4, 8).Value = 0.547381927048329  # H4
$ws.Cells.Item(4, 9).Value = 0.717358098613063  # I4
$ws.Cells.Item(4, 10).Value = 0.000244394207857274  # J4
$ws.Cells.Item(4, 11).Value = 0.0430133805828802  # K4
$ws.Cells.Item(4, 12).Value = 0.00348261746196615  # L4
$ws.Cells.Item(4, 13).Value = 0.117309219771491  # M4
$ws.Cells.Item(4, 14).Value = 0.0543777112482434  # N4
$ws.Cells.Item(4, 15).Value = 0.998961324616607  # O4
$ws.Cells.Item(4, 16).Value = 0.000733182623571821  # P4
$ws.Cells.Item(4, 17).Value = 0  # Q4
$ws.Cells.Item(4, 18).Value = 0.000672084071607503  # R4
$ws.Cells.Item(4, 19).Value = 0.0000610985519643184  # S4
$ws.Cells.Item(4, 20).Value = 0.0000610985519643184  # T4
$ws.Cells.Item(4, 21).Value = 0.000733182623571821  # U4
$ws.Cells.Item(4, 22).Value = 0.0106922465937557  # V4
$ws.Cells.Item(4, 23).Value = 0.905724934319057  # W4
$ws.Cells.Item(4, 24).Value = 0.9991446202725  # X4

# Row 5 ("T" base) - columns B:X
$ws.Cells.Item(5, 2).Value = 0.068308181096108  # B5
$ws.Cells.Item(5, 3).Value = 0.00311602615018024  # C5
$ws.Cells.Item(5, 4).Value = 0.997678255025356  # D5
$ws.Cells.Item(5, 5).Value = 0.000366591311785911  # E5
$ws.Cells.Item(5, 6).Value = 0.000183295655892955  # F5
$ws.Cells.Item(5, 7).Value = 0.0197348322844749  # G5
$ws.Cells.Item(5, 8).Value = 0.384065497647706  # H5
$ws.Cells.Item(5, 9).Value = 0.000427689863750229  # I5
$ws.Cells.Item(5, 10).Value = 0.998900226064642  # J5
$ws.Cells.Item(5, 11).Value = 0.00183295655892955  # K5
$ws.Cells.Item(5, 12).Value = 0.000183295655892955  # L5
$ws.Cells.Item(5, 13).Value = 0.000122197103928637  # M5
$ws.Cells.Item(5, 14).Value = 0.0000610985519643184  # N5
$ws.Cells.Item(5, 15).Value = 0.0000610985519643184  # O5
$ws.Cells.Item(5, 16).Value = 0  # P5
$ws.Cells.Item(5, 17).Value = 0.0799169059693285  # Q5
$ws.Cells.Item(5, 18).Value = 0.0000610985519643184  # R5
$ws.Cells.Item(5, 19).Value = 0.000366591311785911  # S5
$ws.Cells.Item(5, 20).Value = 0.00947027555446936  # T5
$ws.Cells.Item(5, 21).Value = 0.000122197103928637  # U5
$ws.Cells.Item(5, 22).Value = 0.000427689863750229  # V5
$ws.Cells.Item(5, 23).Value = 0.00201625221482251  # W5
$ws.Cells.Item(5, 24).Value = 0  # X5

Write-Output "Updated B2:X5 with new frequency values"
